$wb = $excel.ActiveWorkbook

# Rename the worksheet from "NewRecords" to "Records"
$ws = $wb.ActiveSheet
$ws.Name = "Records"

# Move the active selection on the records sheet from Q3 to Q1
$ws.Range("Q1").Select()
